$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet (fallback to 120 if unavailable)
$lastRow = 120
try {
    $ur = $ws.UsedRange
    $candidate = $ur.Row + $ur.Rows.Count - 1
    if ($candidate -gt 0) {
        $lastRow = $candidate
    }
} catch {
    $lastRow = 120
}

# Column C holds the "Förändrad" (changed) date for each record.
# Update every cell currently holding the old date serial (46061)
# to the new date serial (46062).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
